$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 173. This shifts existing rows 173..283
# down to 174..284, matching the target diff (row 283's old content ends
# up as new row 284; a brand new record is inserted at row 173).
$ws.Rows.Item(173).Insert()

# Populate the newly-inserted row 173 with the new record's data.
$ws.Range("A173").Value = 5
$ws.Range("B173").Value = "Macroferia Regional de Talca"
$ws.Range("C173").Value = "Maule"
$ws.Range("D173").Value = 44596
$ws.Range("E173").Value = 7
$ws.Range("F173").Value = 100112032
$ws.Range("G173").Value = "Zapallo italiano"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 400
$ws.Range("K173").Value = 7000
$ws.Range("L173").Value = 7000
$ws.Range("M173").Value = 7000
$ws.Range("N173").Value = "$/caja 50 unidades"
$ws.Range("O173").Value = "Región del Maule"
$ws.Range("P173").Value = 140
$ws.Range("Q173").Value = 50
$ws.Range("R173").Value = "Hortaliza"
